# Applies the addition of two new dialogue-line rows (rows 8-9) and one
# short-form continuation row (row 10) to the "Сайтер" character sheet,
# matching xl/sharedStrings.xml + xl/worksheets/sheet1.xml target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 8 uses the same bordered style (xf 8 / xf 9) as rows 2-7, so copy
# that formatting (format-only paste) from row 3 before filling values.
# ------------------------------------------------------------------
$ws.Range("A3:E3").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Fill in the new cell values. The order below matches the order new
# strings were typed into the sheet (and therefore the order they were
# appended to the shared string table).
# ------------------------------------------------------------------
$ws.Range("C8").Value = " That [CS:N]Grovyle[CR]...[K] I\'ll rough him up\nwith my claws!"
$ws.Range("A8").Value = "SCRIPT/G01P03A/um1313.ssb"
$ws.Range("D8").Value = " Этот [CS:N]Гровайл[CR]...[K] Отведает моих\nлезвий!"
$ws.Range("E8").Value = " Üóïó [CS:N]Ãñïâàêì[CR]...[K] Ïóâåäàåó íïéö\nìåèâéê!"
$ws.Range("A9").Value = "SCRIPT/T01P02A/um1604.ssb"
$ws.Range("C9").Value = " He must have been a nasty piece\nof work in the future too, eh?"
$ws.Range("C10").Value = " That [CS:N]Grovyle[CR]."
$ws.Range("D9").Value = " Должно быть, он и в будущем\nмного кому насолил, а?"
$ws.Range("D10").Value = " Этот [CS:N]Гровайл[CR]."
$ws.Range("E9").Value = " Äïìçîï áúóû, ïî é â áôäôþåí\níîïãï ëïíô îàòïìéì, à?"
$ws.Range("E10").Value = " Üóïó [CS:N]Ãñïâàêì[CR]."

# Numeric "line number" column (B) - plain numbers, not shared strings.
$ws.Range("B8").Value = 307
$ws.Range("B9").Value = 285
$ws.Range("B10").Value = 288

# ------------------------------------------------------------------
# Row heights for the two full (3-line) dialogue rows; row 10 keeps the
# default row height, same as in the target workbook.
# ------------------------------------------------------------------
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2

# ------------------------------------------------------------------
# Update the view so the newly added rows are visible, mirroring the
# sheetView/selection change recorded in the saved workbook.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D8").Select()
